# The sheet previously held two separate leads in rows 2 and 3 (Alice
# Johnson / Charlie Brown). The source system re-synced: Charlie Brown's
# record was refreshed (new email, cleared country code, lower-cased
# "yes" flag, and a newer call transcript/timestamp) and now occupies
# row 2, while the stale duplicate in row 3 is dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Overwrite row 2 with the refreshed Charlie Brown record ---------------
$ws.Range("A2").Value = 38

# customer_id (B2) and phone_number (D2) are digit-only strings; make sure
# they stay text rather than being auto-detected as numbers, then restore
# the default "Normal" style so no stray number-format is left behind.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "Charlie Brown"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "8107331777"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "vinayak_sharma@technologymindz.com"
$ws.Range("F2").Value = "Charlie is evaluating enterprise-level solutions with a strong emphasis on scalability, integration with his existing ERP system, and compliance with international data protection regulations. He also needs a custom training program for his team."
$ws.Range("G2").Value = "no-answer"

# country_code is now blank.
$ws.Range("H2").Value = ""

$ws.Range("I2").Value = "Real Estate"
$ws.Range("J2").Value = "XYZ Company Ltd."
$ws.Range("K2").Value = "Berlin, Germany"
$ws.Range("L2").Value = "yes"
$ws.Range("M2").Value = "Charlie has shown strong interest in a long-term partnership if the enterprise solution aligns with his company’s compliance and integration needs. He mentioned that decision-making will involve multiple stakeholders, and the procurement cycle might take up to three months. We should prepare detailed documentation, case studies, and a tailored presentation for his board of directors.
[2025-08-25 17:13:21] No summary available. Conversation transcript missing."
$ws.Range("N2").Value = "nan
[2025-08-25 17:13:21] No tasks found for this call."

# --- Drop the now-duplicate row 3 (old Charlie Brown entry) ----------------
$ws.Rows("3:3").Delete()
